$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.736.05"
$ws.Range("E2").Value = '  +0.18%  '
$ws.Range("D3").Value = "'1.601.25"
$ws.Range("E3").Value = '  +0.19%  '
$ws.Range("E4").Value = '  +0.18%  '
$ws.Range("D5").Value = "'211.86"
$ws.Range("E5").Value = '  +0.14%  '
$ws.Range("D6").Value = "'0.513"
$ws.Range("E6").Value = '  +0.02%  '
$ws.Range("E7").Value = '  +0.18%  '
$ws.Range("D8").Value = "'0.0619"
$ws.Range("E8").Value = '  +0.12%  '
$ws.Range("E9").Value = '  +0.15%  '
$ws.Range("D10").Value = "'19.59"
$ws.Range("E10").Value = '  +0.68%  '
$ws.Range("D11").Value = "'0.0848"
$ws.Range("E11").Value = '  +0.79%  '
$ws.Range("D12").Value = "'1.829.31"
$ws.Range("E12").Value = '  +0.36%  '
$ws.Range("D13").Value = "'1.600.35"
$ws.Range("E13").Value = '  +0.12%  '
$ws.Range("E14").Value = '  +0.93%  '
$ws.Range("D15").Value = "'0.525"
$ws.Range("E15").Value = '  +0.37%  '
$ws.Range("D16").Value = "'65.06"
$ws.Range("E16").Value = '  -0.14%  '
$ws.Range("D17").Value = "'0.0₃0740"
$ws.Range("E17").Value = '  -1.43%  '
$ws.Range("E18").Value = '  +0.14%  '
$ws.Range("D19").Value = "'209.28"
$ws.Range("E19").Value = '  -0.20%  '
$ws.Range("D20").Value = "'7.18"
$ws.Range("E20").Value = '  +2.54%  '
$ws.Range("D21").Value = "'4.30"
$ws.Range("E21").Value = '  +0.39%  '
$ws.Range("E22").Value = '  -4.53%  '
$ws.Range("D23").Value = "'9.04"
$ws.Range("E23").Value = '  +0.78%  '
$ws.Range("D24").Value = "'143.48"
$ws.Range("E24").Value = '  +0.32%  '
$ws.Range("E25").Value = '  +0.21%  '
$ws.Range("D26").Value = "'7.12"
$ws.Range("E26").Value = '  +0.18%  '
$ws.Range("E27").Value = '  -0.06%  '
$ws.Range("E28").Value = '  +0.16%  '
$ws.Range("D29").Value = "'0.0507"
$ws.Range("E29").Value = '  -1.62%  '
$ws.Range("E30").Value = '  +0.41%  '
$ws.Range("D31").Value = "'3.27"
$ws.Range("E31").Value = '  +0.45%  '
$ws.Range("E32").Value = '  +0.31%  '
$ws.Range("D33").Value = "'1.282.65"
$ws.Range("E33").Value = '  -0.57%  '
$ws.Range("E34").Value = '  +1.45%  '
$ws.Range("E35").Value = '  +16.25%  '
$ws.Range("E36").Value = '  +0.23%  '
$ws.Range("E37").Value = '  -4.29%  '
$ws.Range("D38").Value = "'0.0169"
$ws.Range("E38").Value = '  -1.17%  '
$ws.Range("D39").Value = "'0.825"
$ws.Range("E39").Value = '  -0.12%  '
$ws.Range("D40").Value = "'5.48"
$ws.Range("E40").Value = '  +0.65%  '
$ws.Range("D41").Value = "'2.19"
$ws.Range("E41").Value = '  +0.19%  '
$ws.Range("E42").Value = '  -0.40%  '
$ws.Range("D43").Value = "'62.65"
$ws.Range("E43").Value = '  -0.96%  '
$ws.Range("D44").Value = "'1.740.44"
$ws.Range("E44").Value = '  +0.32%  '
$ws.Range("D45").Value = "'90.23"
$ws.Range("E45").Value = '  -0.77%  '
$ws.Range("D46").Value = "'1.56"
$ws.Range("E46").Value = '  -0.41%  '
$ws.Range("E47").Value = '  +1.80%  '
$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D48").Value = "'0.0513"
$ws.Range("E48").Value = '  +0.78%  '
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").Value = "'7.56"
$ws.Range("E49").Value = '  +3.26%  '
$ws.Range("B50").Value = 'USDD'
$ws.Range("C50").Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range("D50").Value = "'1.01"
$ws.Range("E50").Value = '  +0.25%  '
$ws.Range("B51").Value = 'Mantle'
$ws.Range("C51").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D51").Value = "'0.399"
$ws.Range("E51").Value = '  +1.56%  '
